$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of form-submission data (row 4), mirroring the existing
# rows (Jan/Faizi..., Ali/Jan...)
$ws.Range("A4").Value = "Anisa"
$ws.Range("B4").Value = "Faizi"
$ws.Range("C4").Value = "anisa@gmail.com"
$ws.Range("D4").Value = 2023439873
$ws.Range("E4").Value = "xyz123"
$ws.Range("F4").Value = "xyz123"

# The email column is a mailto hyperlink, same as rows 2 & 3
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:anisa@gmail.com")

# Re-apply the same visual style as the other email cells (e.g. C3) so the
# new cell matches the existing hyperlink-styled cells exactly
$ws.Range("C4").Style = $ws.Range("C3").Style

# Move / update the active selection like in the saved workbook
$ws.Range("F7").Select()
